$wb = $excel.ActiveWorkbook

$target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e29a24b3f66920984dfe27d5c384205ffb39974/e2e/849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md"
$statusText = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E:F").ColumnWidth = 29.9777047293527

# --- zh-cn sheet ------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C:C").ColumnWidth = 29.9777047293527

$wsZh.Range("I2").Value = "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $target, [System.Reflection.Missing]::Value, "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md", "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md")
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276

$wsZh.Range("J2").Value = "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.6cde7396a0c1c32beacae51b5b1dd5485bd2aa91.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-05 07:08:57"
$wsZh.Range("I:J").ColumnWidth = 40

# --- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C:C").ColumnWidth = 29.9777047293527

$wsDe.Range("I2").Value = "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $target, [System.Reflection.Missing]::Value, "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md", "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md")
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276

$wsDe.Range("J2").Value = "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.6cde7396a0c1c32beacae51b5b1dd5485bd2aa91.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-05 07:09:11"
$wsDe.Range("I:J").ColumnWidth = 40

Write-Output "Report generated for handback"
